$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated "Time Multithreading (s)" values (column D), rows 2-19.
$newValues = @{
    2  = 0.03
    3  = 0.018
    4  = 0.007
    5  = 0.005
    6  = 0.009
    7  = 0.009
    8  = 0.18
    9  = 0.055
    10 = 0.024
    11 = 0.014
    12 = 0.009
    13 = 0.009
    14 = 0.387
    15 = 0.134
    16 = 0.055
    17 = 0.06
    18 = 0.049
    19 = 0.039
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}

# Reflect the author's final active-cell selection.
$ws.Range("D30").Select()
